$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking price cells so they are not
# reinterpreted as numbers (original workbook stores all prices as text).
$textCells = @('D5','D8','D9','D10','D11','D17','D20','D21','D22','D23','D26','D28','D29','D31','D32','D36','D40','D43','D44','D46','D47','D48','D49','D50')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '25.761.63'
$ws.Range('E2').Value = '  -0.38%  '
$ws.Range('D3').Value = '1.633.56'
$ws.Range('E3').Value = '  -0.28%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '215.03'
$ws.Range('E5').Value = '  -0.21%  '
$ws.Range('E6').Value = '  -0.12%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').Value = '0.258'
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '0.0640'
$ws.Range('E9').Value = '  -0.47%  '
$ws.Range('D10').Value = '19.80'
$ws.Range('E10').Value = '  +0.54%  '
$ws.Range('D11').Value = '0.0779'
$ws.Range('E11').Value = '  -0.09%  '
$ws.Range('E12').Value = '  -0.95%  '
$ws.Range('D13').Value = '1.638.81'
$ws.Range('E13').Value = '  +0.05%  '
$ws.Range('D14').Value = '1.860.73'
$ws.Range('E14').Value = '  -0.24%  '
$ws.Range('E15').Value = '  -1.01%  '
$ws.Range('D16').Value = '0.0₃0775'
$ws.Range('E16').Value = '  +1.64%  '
$ws.Range('D17').Value = '63.06'
$ws.Range('E17').Value = '  +0.05%  '
$ws.Range('D18').Value = '25.782.12'
$ws.Range('E18').Value = '  -0.40%  '
$ws.Range('E19').Value = '  -0.11%  '
$ws.Range('D20').Value = '4.43'
$ws.Range('E20').Value = '  +2.46%  '
$ws.Range('D21').Value = '193.94'
$ws.Range('E21').Value = '  -0.27%  '
$ws.Range('D22').Value = '9.95'
$ws.Range('E22').Value = '  +0.61%  '
$ws.Range('D23').Value = '6.13'
$ws.Range('E23').Value = '  +0.57%  '
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('E25').Value = '  -1.10%  '
$ws.Range('D26').Value = '139.60'
$ws.Range('E26').Value = '  -0.32%  '
$ws.Range('E27').Value = '  -4.38%  '
$ws.Range('D28').Value = '6.84'
$ws.Range('E28').Value = '  +0.15%  '
$ws.Range('D29').Value = '15.56'
$ws.Range('E29').Value = '  +0.90%  '
$ws.Range('E30').Value = '  -0.12%  '
$ws.Range('D31').Value = '0.0490'
$ws.Range('E31').Value = '  +0.28%  '
$ws.Range('D32').Value = '3.34'
$ws.Range('E32').Value = '  +1.51%  '
$ws.Range('E33').Value = '  +1.07%  '
$ws.Range('E34').Value = '  +1.22%  '
$ws.Range('E35').Value = '  +0.33%  '
$ws.Range('D36').Value = '0.896'
$ws.Range('E36').Value = '  -0.86%  '
$ws.Range('E37').Value = '  -0.28%  '
$ws.Range('E38').Value = '  +0.17%  '
$ws.Range('D39').Value = '1.106.78'
$ws.Range('E39').Value = '  -1.64%  '
$ws.Range('D40').Value = '0.0156'
$ws.Range('E40').Value = '  +0.12%  '
$ws.Range('E41').Value = '  +0.39%  '
$ws.Range('E42').Value = '  +0.21%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').Value = '99.18'
$ws.Range('E43').Value = '  +1.14%  '
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').Value = '0.801'
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('D45').Value = '0.0₆0109'
$ws.Range('E45').Value = '  -3.82%  '
$ws.Range('D46').Value = '55.13'
$ws.Range('E46').Value = '  -0.52%  '
$ws.Range('D47').Value = '2.47'
$ws.Range('E47').Value = '  +12.02%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '7.71'
$ws.Range('E48').Value = '  -0.43%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').Value = '0.418'
$ws.Range('E49').Value = '  -1.95%  '
$ws.Range('D50').Value = '0.0503'
$ws.Range('E50').Value = '  -0.16%  '
$ws.Range('E51').Value = '  -0.77%  '

# Reset style back to the default (Normal) so no stray number format
# style is left behind on these cells, matching original formatting.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
